$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: Replace the sentence describing sqlite/csv reformatting with
# new text describing converting pmids/dates and unlisting author names.
# ---------------------------------------------------------------------
$old1 = "reads it using the pandas read_csv function. This data is then reformatted for use with sqlite and saved as a new csv file called hiv_csv. "
$new1 = "reads it using the pandas read_csv function, converts pmids to int tye, dates to datetime and unlists the author names because nested lists aren" + [char]0x2019 + "t appropriate in a data frame."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------
# Edit 2: Replace the "sqlalchemy create engine" sentence with the new
# "query the previously created database" phrasing.
# ---------------------------------------------------------------------
$old2 = "Using similar syntax and commands as before, we use a similar create engine function included in sqlalchemy to automatically build a database, specifying "
$new2 = "Using similar syntax and commands as before, we query the previously created database, specifying "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# ---------------------------------------------------------------------
# Edit 3: Move the hidden "_GoBack" bookmark from the end of the
# "...query using the head function. " paragraph to the end of the
# paragraph that now ends with "...a data frame." (i.e. to right after
# the last edit location), matching Word's normal behaviour of tracking
# the most recent edit point.
# ---------------------------------------------------------------------

# Remove the bookmark from its old location (end of the last paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the paragraph that now ends with "...a data frame."
$findRng = $d.Content
$findRng.Find.Execute("a data frame.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $findRng.Paragraphs(1)
$insertPos = $targetPara.Range.End - 1

# NOTE: this COM-interop runtime has a quirk where adding a bookmark with
# a *collapsed* range exactly at a paragraph boundary (End-1 or End)
# silently produces a bogus bookmark elsewhere. Work around it by
# temporarily inserting a marker character, wrapping it with the
# bookmark (a non-collapsed range, which works correctly), and then
# deleting the marker character again -- this correctly leaves behind a
# properly collapsed bookmark at the desired position.
$markerRange = $d.Range($insertPos, $insertPos)
$markerRange.InsertBefore("@")
$wrapRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $wrapRange)
$deleteRange = $d.Range($insertPos, $insertPos + 1)
$deleteRange.Delete()
